# edit.ps1 - applies the ResumeJobHelperTemplate.docx changes described by the diff:
#   1. Lowercase the placeholder tokens {{Frontend}} / {{Backend}} / {{Databases}} / {{Tools}}
#      to {{frontend}} / {{backend}} / {{databases}} / {{tools-cloud}} (run-split preserved)
#   2. Resize the three "job" tables' two columns (7735/1789 or 7645/1879 -> 7015/2509)
#      and switch the right-hand date-range cell justification from "both" to "right"
#   3. Right-justify the third (most recent) job table as a whole
#   4. Collapse the three {{Certification_N}} bullet paragraphs into a single
#      non-list {{certifications}} paragraph

$d = $word.ActiveDocument

function Replace-RunSplit($findText, $xmlFragment) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $target = $d.Range($rng.Start, $rng.End)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $xmlFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

# --- 1. Technical skills placeholders ---------------------------------

Replace-RunSplit "Frontend}}" '<w:p><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:bCs/><w:sz w:val="20"/></w:rPr><w:t>f</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:bCs/><w:sz w:val="20"/></w:rPr><w:t>rontend}}</w:t></w:r></w:p>'

Replace-RunSplit "Backend}}" '<w:p><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:bCs/><w:sz w:val="20"/></w:rPr><w:t>b</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:bCs/><w:sz w:val="20"/></w:rPr><w:t>ackend}}</w:t></w:r></w:p>'

Replace-RunSplit "Databases}}" '<w:p><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:sz w:val="20"/></w:rPr><w:t>d</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:sz w:val="20"/></w:rPr><w:t>atabases}}</w:t></w:r></w:p>'

Replace-RunSplit "Tools}}" '<w:p><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:sz w:val="20"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:sz w:val="20"/></w:rPr><w:t>ools</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:sz w:val="20"/></w:rPr><w:t>-cloud</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:sz w:val="20"/></w:rPr><w:t>}}</w:t></w:r></w:p>'

# --- 2. Resize job-history tables & fix right column justification -----

for ($i = 1; $i -le 3; $i++) {
    $t = $d.Tables($i)
    $t.Columns(1).Width = 350.75   # 7015 dxa
    $t.Columns(2).Width = 125.45   # 2509 dxa
    $secondRowRightCell = $t.Cell(2, 2)
    $secondRowRightCell.Range.Paragraphs(1).Alignment = 2   # wdAlignParagraphRight
}

# --- 3. Right-justify the third (Codewiz) table -------------------------

$t3 = $d.Tables(3)
$t3.Alignment = 2                 # wdAlignParagraphRight  -> tblPr/jc
$t3.Rows(1).Alignment = 2         # trPr/jc on row 1
$t3.Rows(2).Alignment = 2         # trPr/jc on row 2

# --- 4. Collapse the Certification_1/2/3 paragraphs into one -----------

$rngC1 = $d.Content
$rngC1.Find.Execute("{{Certification_1}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $rngC1.Paragraphs(1).Range.Start

$rngC3 = $d.Content
$rngC3.Find.Execute("{{Certification_3}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $rngC3.Paragraphs(1).Range.End

$certTarget = $d.Range($startPos, $endPos)
$certXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="94" w:line="301" w:lineRule="auto"/><w:ind w:left="360"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{{</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>c</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>ertification</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$certTarget.InsertXML($certXml)

Write-Output "done"
